$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.75670566666667
$ws.Range("H2").Value = 59.270117
$ws.Range("I2").Value = 0.05135788836328295
$ws.Range("J2").Value = 0.05135788836328296
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 2222.893969302288
$ws.Range("R2").Value = 20006.04572372059
$ws.Range("S2").Value = 0.01682131582465066
$ws.Range("T2").Value = 0.01682131582465066
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.75670566666667
$ws.Range("H3").Value = 59.270117
$ws.Range("I3").Value = 0.05135788836328295
$ws.Range("J3").Value = 0.05135788836328296
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 2100.42361287084
$ws.Range("R3").Value = 18903.81251583757
$ws.Range("S3").Value = 0.01589454532945806
$ws.Range("T3").Value = 0.01589454532945807
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.75670566666667
$ws.Range("H4").Value = 59.270117
$ws.Range("I4").Value = 0.05135788836328295
$ws.Range("J4").Value = 0.05135788836328296
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 2463.496333510114
$ws.Range("R4").Value = 22171.46700159102
$ws.Range("S4").Value = 0.01864202720917423
$ws.Range("T4").Value = 0.01864202720917423
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 329.6209206666667
$ws.Range("H5").Value = 988.862762
$ws.Range("I5").Value = 0.8568551220744788
$ws.Range("J5").Value = 0.8568551220744789
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 37086.76785836957
$ws.Range("R5").Value = 333780.9107253261
$ws.Range("S5").Value = 0.2806468701055265
$ws.Range("T5").Value = 0.2806468701055265
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 329.6209206666667
$ws.Range("H6").Value = 988.862762
$ws.Range("I6").Value = 0.8568551220744788
$ws.Range("J6").Value = 0.8568551220744789
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 35043.47216310503
$ws.Range("R6").Value = 315391.2494679453
$ws.Range("S6").Value = 0.2651846291314408
$ws.Range("T6").Value = 0.2651846291314408
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 329.6209206666667
$ws.Range("H7").Value = 988.862762
$ws.Range("I7").Value = 0.8568551220744788
$ws.Range("J7").Value = 0.8568551220744789
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 41100.97823042401
$ws.Range("R7").Value = 369908.8040738161
$ws.Range("S7").Value = 0.3110236228375115
$ws.Range("T7").Value = 0.3110236228375116
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 35.30925033333333
$ws.Range("H8").Value = 105.927751
$ws.Range("I8").Value = 0.09178698956223814
$ws.Range("J8").Value = 0.09178698956223817
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 3972.763523980464
$ws.Range("R8").Value = 35754.87171582417
$ws.Range("S8").Value = 0.03006311180667915
$ws.Range("T8").Value = 0.03006311180667916
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 35.30925033333333
$ws.Range("H9").Value = 105.927751
$ws.Range("I9").Value = 0.09178698956223814
$ws.Range("J9").Value = 0.09178698956223817
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 3753.884094048655
$ws.Range("R9").Value = 33784.9568464379
$ws.Range("S9").Value = 0.02840678448326745
$ws.Range("T9").Value = 0.02840678448326747
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 35.30925033333333
$ws.Range("H10").Value = 105.927751
$ws.Range("I10").Value = 0.09178698956223814
$ws.Range("J10").Value = 0.09178698956223817
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 4402.768872642385
$ws.Range("R10").Value = 39624.91985378147
$ws.Range("S10").Value = 0.03331709327229154
$ws.Range("T10").Value = 0.03331709327229155